$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update issue statuses (rows 2-6 in the pre-sort layout):
#  Issue_031 "In work" -> "Solved"
#  Issue_034 "In work" -> "Solved"
#  Issue_032 "Open"    -> "Solved"
#  Issue_010 "Open"    -> "In Work"
#  Issue_033 "Open"    -> "Solved"
$ws.Range("D2").Value = "Solved"
$ws.Range("D3").Value = "Solved"
$ws.Range("D4").Value = "Solved"
$ws.Range("D5").Value = "In Work"
$ws.Range("D6").Value = "Solved"

# Re-sort the issue list the same way the sheet's existing sortState does:
# Status ascending, then Priority descending, then ID ascending.
$rng = $ws.Range("A1:F35")
$key1 = $ws.Range("D1")
$key2 = $ws.Range("C1")
$key3 = $ws.Range("A1")
$rng.Sort($key1, 1, $key2, $null, 2, $key3, 1, 1)

# Move the active selection
$ws.Range("D3").Select()

$wb.Save()
